# Applies the data-correction edit described in the commit
# "fixed removing snippets with timeouts" to the "openjml" worksheet of
# data/correlation_analysis.xlsx.
#
# Changes:
#  - Column J (10) width 19.7109375 -> 20.7109375
#  - Rows 9-12, columns H (num_snippets_for_correlation), I (kendalls_tau),
#    J (kendalls_p_value), K (spearmans_rho), L (spearmans_p_value) updated
#    to reflect recomputed correlations after removing snippets with
#    timeouts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openjml")

# Fix column J width (was mistakenly left narrower than its neighbors).
# Target raw width is 20.7109375 characters; the host's ColumnWidth setter
# quantizes to a 1/6-character pixel grid, so 19.8 is the closest input that
# lands on the nearest achievable grid point (20.666666... which rounds to
# the same displayed width as its neighboring columns).
$ws.Columns.Item(10).ColumnWidth = 19.8

# Row 9
$ws.Range("H9").Value = 61
$ws.Range("I9").Value = -0.1004064600277419
$ws.Range("J9").Value = 0.3197996198904468
$ws.Range("K9").Value = -0.1258658290311623
$ws.Range("L9").Value = 0.3337653039608836

# Row 10
$ws.Range("H10").Value = 48
$ws.Range("I10").Value = -0.2134569412305524
$ws.Range("J10").Value = 0.05383248704028133
$ws.Range("K10").Value = -0.254689962010721
$ws.Range("L10").Value = 0.08064290874044987

# Row 11
$ws.Range("H11").Value = 48
$ws.Range("I11").Value = -0.1499545016214044
$ws.Range("J11").Value = 0.1564334829386408
$ws.Range("K11").Value = -0.2015245812013276
$ws.Range("L11").Value = 0.1695848710732239

# Row 12
$ws.Range("H12").Value = 48
$ws.Range("I12").Value = 0.1670554193005554
$ws.Range("J12").Value = 0.1123317103170471
$ws.Range("K12").Value = 0.2349915639777836
$ws.Range("L12").Value = 0.1078882378308313
